$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 1.68
$ws.Range("H2").Value = 5.9
$ws.Range("N2").Value = 4.1
$ws.Range("Q2").Value = 1.88
$ws.Range("S2").Value = 3.25
$ws.Range("AA2").Value = 160
$ws.Range("G5").Value = 7
$ws.Range("H5").Value = 1.55
$ws.Range("K5").Value = 4.9
$ws.Range("Q5").Value = 1.71
$ws.Range("S5").Value = 2.78
$ws.Range("V5").Value = 2.56
$ws.Range("W5").Value = 1.17
$ws.Range("AD5").Value = 11
$ws.Range("AF5").Value = 55
$ws.Range("AG5").Value = 25
$ws.Range("AH5").Value = 990
$ws.Range("AK5").Value = 95
$ws.Range("Q6").Value = 1.64
$ws.Range("W6").Value = 3.15
$ws.Range("F7").Value = 1.72
$ws.Range("G7").Value = 1.88
$ws.Range("J7").Value = 3.55
$ws.Range("Q7").Value = 1.97
$ws.Range("W7").Value = 2.12
$ws.Range("F8").Value = 1.41
$ws.Range("G8").Value = 1.42
$ws.Range("P8").Value = 2.04
$ws.Range("Q8").Value = 1.94
$ws.Range("V8").Value = 1.1
$ws.Range("AL8").Value = 44
$ws.Range("S9").Value = 3.8
$ws.Range("X9").Value = 14.5
$ws.Range("G11").Value = 7.4
$ws.Range("J11").Value = 4.1
$ws.Range("O11").Value = 1.38
$ws.Range("G12").Value = 2.88
$ws.Range("H12").Value = 2.86
$ws.Range("I12").Value = 2.9
$ws.Range("J12").Value = 3.25
$ws.Range("P12").Value = 1.71
$ws.Range("Q12").Value = 2.34
$ws.Range("V12").Value = 1.52
$ws.Range("Y12").Value = 9.6
$ws.Range("AB12").Value = 9.800000000000001
$ws.Range("AL12").Value = 60
$ws.Range("AN12").Value = 36
$ws.Range("T13").Value = 1.99
$ws.Range("Q14").Value = 2.48
$ws.Range("T14").Value = 2.02
$ws.Range("P15").Value = 2.06
$ws.Range("R15").Value = 1.34
$ws.Range("S15").Value = 2.68
$ws.Range("T15").Value = 1.5
$ws.Range("U15").Value = 1.84
$ws.Range("AD15").Value = 23
$ws.Range("AN15").Value = 18.5
$ws.Range("F16").Value = 2.3
$ws.Range("G16").Value = 2.32
$ws.Range("P16").Value = 1.92
$ws.Range("V16").Value = 1.38
$ws.Range("AE16").Value = 40
$ws.Range("H17").Value = 5.4
$ws.Range("I17").Value = 5.7
$ws.Range("M18").Value = 1.08
$ws.Range("Q18").Value = 2.12
$ws.Range("AM18").Value = 140
$ws.Range("I19").Value = 2.12
$ws.Range("X19").Value = 16
$ws.Range("N20").Value = 7.2
$ws.Range("O20").Value = 1.13
$ws.Range("R20").Value = 1.86
$ws.Range("S20").Value = 2
$ws.Range("T20").Value = 2.28
$ws.Range("U20").Value = 1.64
$ws.Range("X20").Value = 990
$ws.Range("Y20").Value = 990
$ws.Range("Z20").Value = 290
$ws.Range("AB20").Value = 13
$ws.Range("AC20").Value = 23
$ws.Range("AD20").Value = 990
$ws.Range("AE20").Value = 470
$ws.Range("AF20").Value = 9.199999999999999
$ws.Range("AG20").Value = 14.5
$ws.Range("AH20").Value = 990
$ws.Range("AI20").Value = 310
$ws.Range("AJ20").Value = 9
$ws.Range("AK20").Value = 15.5
$ws.Range("AL20").Value = 980
$ws.Range("AM20").Value = 280
$ws.Range("AN20").Value = 3.05
$ws.Range("J21").Value = 1.03
$ws.Range("L21").Value = 1.01
$ws.Range("M21").Value = 1.01
$ws.Range("N21").Value = 1.28
$ws.Range("O21").Value = 1.28
$ws.Range("P21").Value = 1.28
$ws.Range("Q21").Value = 1.28
$ws.Range("R21").Value = 1.18
$ws.Range("S21").Value = 1.28
$ws.Range("T21").Value = 1.04
$ws.Range("U21").Value = 1.04
$ws.Range("V21").Value = 1.01
$ws.Range("W21").Value = 1.01
$ws.Range("X21").Value = 990
$ws.Range("Y21").Value = 990
$ws.Range("Z21").Value = 1000
$ws.Range("AA21").Value = 1000
$ws.Range("AB21").Value = 990
$ws.Range("AC21").Value = 990
$ws.Range("AD21").Value = 990
$ws.Range("AE21").Value = 1000
$ws.Range("AF21").Value = 1000
$ws.Range("AG21").Value = 990
$ws.Range("AH21").Value = 990
$ws.Range("AI21").Value = 1000
$ws.Range("AJ21").Value = 1000
$ws.Range("AK21").Value = 1000
$ws.Range("AL21").Value = 1000
$ws.Range("AM21").Value = 1000
$ws.Range("AN21").Value = 1000
$ws.Range("AO21").Value = 1000
